# web 122 / quiz 02
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new quiz column
$ws.Range("C1").Value = "Q02"

# Quiz 02 raw scores -> scaled to /10, entered as formulas like the other quiz columns
$ws.Range("C2").Formula  = "=(21/25)*10"
$ws.Range("C3").Formula  = "=(18/25)*10"
$ws.Range("C4").Formula  = "=(21/25)*10"
$ws.Range("C5").Formula  = "=(22/25)*10"
$ws.Range("C6").Formula  = "=(23/25)*10"
$ws.Range("C7").Formula  = "=(24/25)*10"
$ws.Range("C8").Formula  = "=(22/25)*10"
$ws.Range("C9").Formula  = "=(15/25)*10"
$ws.Range("C10").Formula = "=(21/25)*10"
$ws.Range("C11").Formula = "=(15/25)*10"
$ws.Range("C12").Formula = "=(21/25)*10"

# Move the active selection to reflect where the author ended up
$ws.Range("C13").Select()
